$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells hold plain-text price strings (e.g. "213.29"). Assigning
# .Value directly lets the engine auto-detect these as numbers, which would
# change the cell's stored type. Forcing NumberFormat to "@" ("Text") right
# before the write keeps the value a literal string; immediately borrowing
# the NumberFormat/Style from the same row's column-C cell (plain, unstyled)
# afterwards avoids leaving a stray style index attached to the cell.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.659.61'
$ws.Range('D2').NumberFormat = $ws.Range('C2').NumberFormat
$ws.Range('D2').Style = $ws.Range('C2').Style
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.637.14'
$ws.Range('D3').NumberFormat = $ws.Range('C3').NumberFormat
$ws.Range('D3').Style = $ws.Range('C3').Style
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.29'
$ws.Range('D5').NumberFormat = $ws.Range('C5').NumberFormat
$ws.Range('D5').Style = $ws.Range('C5').Style
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +2.38%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0624'
$ws.Range('D9').NumberFormat = $ws.Range('C9').NumberFormat
$ws.Range('D9').Style = $ws.Range('C9').Style
$ws.Range('E9').Value = '  +1.43%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.16'
$ws.Range('D10').NumberFormat = $ws.Range('C10').NumberFormat
$ws.Range('D10').Style = $ws.Range('C10').Style
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +3.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.864.82'
$ws.Range('D12').NumberFormat = $ws.Range('C12').NumberFormat
$ws.Range('D12').Style = $ws.Range('C12').Style
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.634.37'
$ws.Range('D13').NumberFormat = $ws.Range('C13').NumberFormat
$ws.Range('D13').Style = $ws.Range('C13').Style
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '26.664.60'
$ws.Range('D16').NumberFormat = $ws.Range('C16').NumberFormat
$ws.Range('D16').Style = $ws.Range('C16').Style
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.33'
$ws.Range('D17').NumberFormat = $ws.Range('C17').NumberFormat
$ws.Range('D17').Style = $ws.Range('C17').Style
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0₃0744'
$ws.Range('D18').NumberFormat = $ws.Range('C18').NumberFormat
$ws.Range('D18').Style = $ws.Range('C18').Style
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '218.53'
$ws.Range('D19').NumberFormat = $ws.Range('C19').NumberFormat
$ws.Range('D19').Style = $ws.Range('C19').Style
$ws.Range('E19').Value = '  +7.73%  '
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.22'
$ws.Range('D23').NumberFormat = $ws.Range('C23').NumberFormat
$ws.Range('D23').Style = $ws.Range('C23').Style
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '149.13'
$ws.Range('D25').NumberFormat = $ws.Range('C25').NumberFormat
$ws.Range('D25').Style = $ws.Range('C25').Style
$ws.Range('E25').Value = '  +4.06%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.85'
$ws.Range('D28').NumberFormat = $ws.Range('C28').NumberFormat
$ws.Range('D28').Style = $ws.Range('C28').Style
$ws.Range('E28').Value = '  +4.37%  '
$ws.Range('E29').Value = '  +1.48%  '
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.31'
$ws.Range('D32').NumberFormat = $ws.Range('C32').NumberFormat
$ws.Range('D32').Style = $ws.Range('C32').Style
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.196.97'
$ws.Range('D36').NumberFormat = $ws.Range('C36').NumberFormat
$ws.Range('D36').Style = $ws.Range('C36').Style
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('E37').Value = '  +5.37%  '
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.506'
$ws.Range('D40').NumberFormat = $ws.Range('C40').NumberFormat
$ws.Range('D40').Style = $ws.Range('C40').Style
$ws.Range('E40').Value = '  +2.02%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.41'
$ws.Range('D42').NumberFormat = $ws.Range('C42').NumberFormat
$ws.Range('D42').Style = $ws.Range('C42').Style
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.773.15'
$ws.Range('D44').NumberFormat = $ws.Range('C44').NumberFormat
$ws.Range('D44').Style = $ws.Range('C44').Style
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '92.27'
$ws.Range('D45').NumberFormat = $ws.Range('C45').NumberFormat
$ws.Range('D45').Style = $ws.Range('C45').Style
$ws.Range('E45').Value = '  -1.40%  '
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('E47').Value = '  +1.53%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0513'
$ws.Range('D48').NumberFormat = $ws.Range('C48').NumberFormat
$ws.Range('D48').Style = $ws.Range('C48').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.63'
$ws.Range('D49').NumberFormat = $ws.Range('C49').NumberFormat
$ws.Range('D49').Style = $ws.Range('C49').Style
$ws.Range('E49').Value = '  +4.68%  '
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('E51').Value = '  +0.09%  '
